$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Terminal La Palmera de La Serena"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 45043
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 100112030
$ws.Range("G9").Value = "Poroto granado"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 560
$ws.Range("K9").Value = 35000
$ws.Range("L9").Value = 36000
$ws.Range("M9").Value = 35500
$ws.Range("N9").Value = "$/malla 25 kilos"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 1420
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
